$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows appended to the log.
# Each entry: Date (text), Number (text, phone-like), Name (text), Balance (numeric)
$data = @(
    ,@('2026-02-13 18:58:04', '237650651854', 'SYLVIANE LAURE SIMO TAGNE ETS MOBILE FINANCIAL SERVICES MFS', 67307)
    ,@('2026-02-13 16:14:55', '237650988697', 'FABASSOU GASPARD ETS MOBILE FINANCIAL SERVICES MFS', 324536)
    ,@('2026-02-13 16:02:17', '237651809692', 'AUREL KAMTA TEUMEN', 224523)
    ,@('2026-02-13 15:05:48', '237652579681', 'AGRIPINE LAURE NOUMUNA', 74577)
    ,@('2026-02-13 15:00:11', '237652724076', 'Nzepena Tambo Jeanine ETS MOBILE FINANCIAL SERVICES MFS', 78924)
    ,@('2026-02-12 05:02:51', '237652958984', 'DEUCHE FEUZE CHARNELLE ETS MOBILE FINANCIAL SERVICES MFS', 324514)
    ,@('2026-02-13 16:47:35', '237653239050', 'ARMAND KAMTA KUITSA', 157987)
    ,@('2026-02-13 14:52:14', '237653423810', 'TEFFO ALAIN NOEL ETS MOBILE FINANCIAL SERVICES MFS', 46342)
    ,@('2026-02-13 16:41:36', '237653472620', 'LANKEU ALAIN JULIOT ETS MOBILE FINANCIAL SERVICES MFS', 314642)
    ,@('2026-02-13 16:50:30', '237654117741', 'DONGMEZA TSAGUE STEPHANE ETS MOBILE FINANCIAL SERVICES MFS', 78187)
    ,@('2026-02-13 15:31:33', '237654134124', 'HERVE ANTOINE EDZILI TCHAMDA', 10490)
    ,@('2026-02-13 15:34:58', '237670358585', 'EVAGLE EMMANUEL -MFS', 203366)
    ,@('2026-02-13 14:49:39', '237670806337', 'MFS AM ANGE RAPHAËL', 0)
    ,@('2026-02-13 16:03:36', '237670809743', 'EBONGUE NYAMSI YVAN LANDRY ETS MOBILE FINANCIAL SERVICES MFS', 103009)
    ,@('2026-02-13 13:56:29', '237670984670', 'LIONELFONE MEDJEUFEUSSI', 7154)
    ,@('2026-02-13 12:56:13', '237671364668', 'ODETTE PELANTITCHOUWE', 93543)
    ,@('2026-02-13 16:53:47', '237671369692', 'EDITH MELI TEMZENETCHOU', 23751)
    ,@('2026-02-13 14:08:25', '237674272162', 'ETS MFS LTDLA CB0114 TCHEUTCHOUA TCHAWE PAUL EUGENE', 296526)
    ,@('2026-02-13 14:23:12', '237674431977', 'TCHETGNIA TOBET RUTH TOP MOBIL', 1088389)
    ,@('2026-02-13 08:57:47', '237674913584', 'MBAKOP DJOUMOU DARC TOP MOBIL', 155195)
    ,@('2026-02-13 17:15:07', '237675432472', 'JUDITH LAURE LIENOU', 3172)
    ,@('2026-02-13 11:12:21', '237675450046', 'BOKO LYNDA JOELLE ETS MOBILE FINANCIAL SERVICES MFS', 137029)
    ,@('2026-02-13 16:21:04', '237675641522', 'LOUISE ODILE DJOUALA', 223390)
    ,@('2026-02-13 16:10:50', '237676270842', 'KAOU POINT SARL SODISERV', 29264)
    ,@('2026-02-13 15:14:45', '237677525770', 'TOP MOBILE LTDLA CB05277 DOMOU ANTOINETTE', 37658)
    ,@('2026-02-13 15:22:58', '237679071551', 'BEGHEH PAMELA MBEBUI ETS MOBILE FINANCIAL SERVICES MFS', 214488)
    ,@('2026-02-13 15:24:59', '237680472978', 'NICOLAS BIKECK', 204326)
    ,@('2026-02-13 12:27:56', '237681102046', 'FRANCOISE NKENFACK NKENGMO', 379934)
    ,@('2026-02-13 14:29:58', '237682368679', 'MFS SIM AA 2', 193111)
    ,@('2026-02-13 16:36:56', '237683360459', 'LUCIE MAJOLIE LELE NKANKEU', 376)
    ,@('2026-02-13 14:19:03', '237652899422', 'CATHERINE CHAMBERLINE NZEPANG', 235473)
    ,@('2026-02-13 13:59:13', '237670904526', 'MFS SIM PROVISOIRE 20', 0)
    ,@('2026-02-13 12:59:38', '237671105116', 'MFS  AM FACE HÔPITAL GÉNÉRAL', 22)
    ,@('2026-02-13 19:14:11', '237672916354', 'MAMADOU DIAN BAH LENA GLOBAL', 16050)
    ,@('2026-02-13 04:18:15', '237672920086', 'NAMY NGOKO CLARISSE ROSE VERTINE KAMILAH CONNECTION', 7230)
    ,@('2026-02-13 14:59:19', '237674000053', 'YVETTE LAURE NGANDO SIMO EPSE DIMO', 172010)
    ,@('2026-02-13 15:25:27', '237674841555', 'BEATRICE TCHAMTIEU EPSE NGAMENI', 202341)
    ,@('2026-02-13 16:45:48', '237674899678', 'VIVIANE MADJUIMEKEM FOMEKONG', 117534)
    ,@('2026-02-13 15:43:51', '237676439452', 'SAGNOU BRINDA JOSELINE _DIGITAL BUSINESS SARL', 95871)
    ,@('2026-02-13 17:29:54', '237676695935', 'MFS MAKEPE MATURITE', 326493)
    ,@('2026-02-13 15:29:59', '237677745809', 'JOSEPHINE CLAIRE NGUENKAM KADJI', 297533)
)

$startRow = 183
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Column A: Date, stored as plain text (matches source data)
    $ws.Cells.Item($r, 1).Value = $row[0]

    # Column B: Number string (e.g. phone number) -- a leading quote
    # forces Excel to store it as text instead of coercing to a number.
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]

    # Column C: Name, plain text
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Column D: Balance, numeric
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# The quote-prefix above also stamps a quotePrefix style flag on the
# cells; clear it in one pass so the new rows keep the default styling
# used by the rest of the data (no explicit style), same as the source.
$ws.Range("B183:B223").ClearFormats()

